$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name (C3)
$ws.Range("C3").Value = "Apurba Khan"

# Row 7 - __init__ / Attributes set to input values.
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "client_number: 1999`nfirst_name: ""Apurba""`nlast_name: ""Khan""`nemail_address: ""mkhan14@rrc.ca"""
$ws.Range("G7").Value = "Object created with expected atttribute value based on method inputs."

# Row 8 - invalid client number
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "client_number: None`nfirst_name: ""Apurba""`nlast_name: ""Khan""`nemail_address: ""mkhan14@rrc.ca"""
$ws.Range("G8").Value = "ValueError(""client_number must be numeric."")"

# Row 9 - blank first_name
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "client_number: 1999`nfirst_name: "" ""`nlast_name: ""Khan""`nemail_address: ""mkhan14@rrc.ca"""
$ws.Range("G9").Value = "ValueError(""first_name cannot be blank."")"

# Row 10 - blank last_name
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "client_number: 1999`nfirst_name: ""Apurba""`nlast_name: "" ""`nemail_address: ""mkhan14@rrc.ca"""
$ws.Range("G10").Value = "ValueError(""last_name cannot be blank."")"

# Row 11 - invalid email
$ws.Range("E11").Value = "None"
$ws.Range("F11").Value = "client_number: 1999`nfirst_name: ""Apurba""`nlast_name: ""Khan""`nemail_address: ""mkhan14"""
$ws.Range("G11").Value = "EmailNotValidError: ""email@pixell-river.com"""

# Row 12 - client_number getter
$ws.Range("E12").Value = "Client(1999,`n""Apurba"",`n""Khan"",`n""mkhan14@rrc.ca"")"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "client._Client__client_number`n= 1999"

# Row 13 - first_name getter
$ws.Range("E13").Value = "Client(1999,`n""Apurba"",`n""Khan"",`n""mkhan14@rrc.ca"")"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "client._Client__first_name`n= ""Apurba"""

# Row 14 - last_name getter
$ws.Range("E14").Value = "Client(1999,`n""Apurba"",`n""Khan"",`n""mkhan14@rrc.ca"")"
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = "client._Client__last_name`n= ""Khan"""

# Row 15 - email_address getter
$ws.Range("E15").Value = "Client(1999,`n""Apurba"",`n""Khan"",`n""mkhan14@rrc.ca"")"
$ws.Range("F15").Value = "None"
$ws.Range("G15").Value = "client._Client__email_address`n= ""mkhan14@rrc.ca"""

# Row 16 - __str__
$ws.Range("E16").Value = "Client(1999,`n""Apurba"",`n""Khan"",`n""mkhan14@rrc.ca"")"
$ws.Range("F16").Value = "None"
$ws.Range("G16").Value = "Khan, Apurba [1999] - mkhan14@rrc.ca"
